$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated file names for existing rows 2-25 (column B)
$updates = @{
    2  = "n1_IMG_3174.jpeg"
    3  = "n2_IMG_3174HorFlip.jpeg"
    4  = "n3_IMG_3174HorVertFlip.jpeg"
    5  = "n4_IMG_3174VertFlip.jpeg"
    6  = "n5_IMG_3175.jpeg"
    7  = "n6_IMG_3175HorFlip.jpeg"
    8  = "n7_IMG_3175HorVertFlip.jpeg"
    9  = "n8_IMG_3175VertFlip.jpeg"
    10 = "n9_IMG_3176.jpeg"
    11 = "n10_IMG_3176HorFlip.jpeg"
    12 = "n11_IMG_3176HorVertFlip.jpeg"
    13 = "n12_IMG_3176VertFlip.jpeg"
    14 = "n13_IMG_3177.jpeg"
    15 = "n14_IMG_3177HorFlip.jpeg"
    16 = "n15_IMG_3177HorVertFlip.jpeg"
    17 = "n16_IMG_3177VertFlip.jpeg"
    18 = "n17_IMG_3178.jpeg"
    19 = "n18_IMG_3178HorFlip.jpeg"
    20 = "n19_IMG_3178HorVertFlip.jpeg"
    21 = "n20_IMG_3178VertFlip.jpeg"
    22 = "n21_IMG_3179.jpeg"
    23 = "n22_IMG_3179HorFlip.jpeg"
    24 = "n23_IMG_3179HorVertFlip.jpeg"
    25 = "n24_IMG_3179VertFlip.jpeg"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}

# New rows 26-29
$newRows = @(
    @("n25", "n25_IMG_3180.jpeg", "'True", "no_meltpatch", "negative"),
    @("n26", "n26_IMG_3180HorFlip.jpeg", "'True", "no_meltpatch", "negative"),
    @("n27", "n27_IMG_3180HorVertFlip.jpeg", "'True", "no_meltpatch", "negative"),
    @("n28", "n28_IMG_3180VertFlip.jpeg", "'True", "no_meltpatch", "negative")
)

$rowIndex = 26
foreach ($rowData in $newRows) {
    $ws.Cells.Item($rowIndex, 1).Value = $rowData[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rowData[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rowData[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rowData[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rowData[4]
    $rowIndex++
}

# Update selection to match new range
$ws.Range("A2:E25").Select()
